$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54 (shifts existing rows 54-74 down to 55-75)
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new weekly data entry
$ws.Cells.Item(54, 1).Value = 7
$ws.Cells.Item(54, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(54, 3).Value = "Ñuble"
$ws.Cells.Item(54, 4).Value = 44924
$ws.Cells.Item(54, 5).Value = 16
$ws.Cells.Item(54, 6).Value = 100112026
$ws.Cells.Item(54, 7).Value = "Haba"
$ws.Cells.Item(54, 8).Value = "Sin especificar"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 60
$ws.Cells.Item(54, 11).Value = 15000
$ws.Cells.Item(54, 12).Value = 15000
$ws.Cells.Item(54, 13).Value = 15000
$ws.Cells.Item(54, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(54, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(54, 16).Value = 600
$ws.Cells.Item(54, 17).Value = 25
$ws.Cells.Item(54, 18).Value = "Hortaliza"
